$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (dLacI): tighten bmin/bmax ---
$ws.Cells.Item(4,2).Value = 0.0001    # B4: 1E-3 -> 1E-4
$ws.Cells.Item(4,3).Value = 1         # C4: 10   -> 1

# --- Row 5 (dCit): tighten bmin/bmax, now estimated ---
$ws.Cells.Item(5,2).Value = 0.0001    # B5: 7E-3 -> 1E-4
$ws.Cells.Item(5,3).Value = 1         # C5: 8E-3 -> 1
$ws.Cells.Item(5,6).Value = "yes"     # F5: no   -> yes

# --- Row 9 (LacI_rep_Cit): tighten bmin ---
$ws.Cells.Item(9,2).Value = 0.00001   # B9: 1E-3 -> 1E-5

# --- Row 10 (LacI_rep_Cit_W220F): tighten bmin ---
$ws.Cells.Item(10,2).Value = 0.00001  # B10: 1E-3 -> 1E-5

# --- New parameter row 16: P_4Lacn_LacI ---
$ws.Cells.Item(16,1).Value = "P_4Lacn_LacI"
$ws.Cells.Item(16,2).Value = 0.01
$ws.Cells.Item(16,3).Value = 100
$ws.Cells.Item(16,4).Value = 1
$ws.Cells.Item(16,5).Value = 98
$ws.Cells.Item(16,6).Value = "yes"
$ws.Cells.Item(16,7).Value = "k_{LacI_W220F_Q60G_T167A}"

# --- New parameter row 17: P_4Lacn_LacI_L ---
$ws.Cells.Item(17,1).Value = "P_4Lacn_LacI_L"
$ws.Cells.Item(17,2).Value = 0.00001
$ws.Cells.Item(17,3).Value = 0.01
$ws.Cells.Item(17,4).Value = 1
$ws.Cells.Item(17,5).Value = 0.0003
$ws.Cells.Item(17,6).Value = "yes"
$ws.Cells.Item(17,7).Value = "kL_W220F_Q60G_T167A"

# --- New parameter row 18: LacI_rep ---
$ws.Cells.Item(18,1).Value = "LacI_rep"
$ws.Cells.Item(18,2).Value = 0.00001
$ws.Cells.Item(18,3).Value = 100
$ws.Cells.Item(18,4).Value = 1
$ws.Cells.Item(18,5).Value = 1
$ws.Cells.Item(18,6).Value = "yes"
$ws.Cells.Item(18,7).Value = "theta_{LacI_W220F_Q60G_T167A}"

# --- Move/save the active selection to H21, matching the saved view state ---
$null = $ws.Range("H21").Select()
